# Apply the "Updated cryptos list" data refresh to Sheet1.
# All target cells in columns B:E are text (inline strings), so we just
# assign string values directly to the corresponding Range cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "91.601.21"
$ws.Range("E2").Value = "  +1.15%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.097.64"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.40%  "

# Row 5 - Solana
$ws.Range("D5").Value = "239.67"
$ws.Range("E5").Value = "  -0.82%  "

# Row 6 - BNB
$ws.Range("D6").Value = "615.17"
$ws.Range("E6").Value = "  -0.85%  "

# Row 7 - XRP
$ws.Range("D7").Value = "1.09"
$ws.Range("E7").Value = "  -5.43%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "0.392"
$ws.Range("E8").Value = "  +7.17%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.01%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.093.45"
$ws.Range("E10").Value = "  -0.23%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -0.96%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.74%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +1.27%  "

# Rows 14/15 swap: WrappedBTC <-> Toncoin with refreshed values
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").Value = "5.50"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "91.552.81"
$ws.Range("E15").Value = "  +1.42%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "33.92"
$ws.Range("E16").Value = "  -3.08%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.668.11"
$ws.Range("E17").Value = "  -0.05%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.098.45"
$ws.Range("E18").Value = "  +0.51%  "

# Row 19 - SuiNetwork
$ws.Range("E19").Value = "  -3.02%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "14.74"
$ws.Range("E20").Value = "  +0.56%  "

# Row 21 - Polkadot
$ws.Range("D21").Value = "5.79"
$ws.Range("E21").Value = "  -0.11%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "9.29"
$ws.Range("E22").Value = "  +2.59%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "443.34"
$ws.Range("E23").Value = "  +0.85%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  -5.22%  "

# Row 25 - NEARProtocol
$ws.Range("D25").Value = "5.60"
$ws.Range("E25").Value = "  -1.04%  "

# Rows 26/27 swap: Litecoin <-> Aptos with refreshed values
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "11.62"
$ws.Range("E26").Value = "  -2.25%  "

$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "80.78"
$ws.Range("E27").Value = "  -10.13%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "3.261.52"
$ws.Range("E28").Value = "  +0.51%  "

# Row 29 - Dai
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.136"
$ws.Range("E30").Value = "  +14.14%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "0.226"
$ws.Range("E31").Value = "  -7.85%  "

# Row 32 - Cronos
$ws.Range("D32").Value = "0.166"
$ws.Range("E32").Value = "  -5.73%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "9.22"
$ws.Range("E33").Value = "  -0.16%  "

# Row 34 - Binance-PegBSC-USD
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  +3.06%  "

# Row 35 - Kaspa
$ws.Range("D35").Value = "0.167"
$ws.Range("E35").Value = "  +0.49%  "

# Row 36 - RenderToken
$ws.Range("D36").Value = "7.95"
$ws.Range("E36").Value = "  +1.04%  "

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "26.11"
$ws.Range("E37").Value = "  -0.67%  "

# Row 38 - MantraDAO
$ws.Range("D38").Value = "4.13"
$ws.Range("E38").Value = "  -3.74%  "

# Row 39 - PancakeSwap
$ws.Range("E39").Value = "  +0.47%  "

# Rows 40/41 swap: Fetch.AI <-> Bittensor with refreshed values
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "479.52"
$ws.Range("E40").Value = "  -2.27%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "1.30"
$ws.Range("E41").Value = "  -0.41%  "

# Row 42 - PolygonEcosystemToken
$ws.Range("D42").Value = "0.430"
$ws.Range("E42").Value = "  +2.25%  "

# Row 43 - dogwifhat
$ws.Range("D43").Value = "3.39"
$ws.Range("E43").Value = "  -5.45%  "

# Row 44 - WhiteBITCoin
$ws.Range("D44").Value = "22.16"
$ws.Range("E44").Value = "  -0.05%  "

# Row 46 - Monero
$ws.Range("D46").Value = "158.85"
$ws.Range("E46").Value = "  +2.40%  "

# Row 47 - Stacks
$ws.Range("D47").Value = "1.91"
$ws.Range("E47").Value = "  -0.86%  "

# Row 48 - ARBITRUM
$ws.Range("D48").Value = "0.694"
$ws.Range("E48").Value = "  +0.72%  "

# Rows 49/50 swap: VeChain <-> ImmutableX with refreshed values
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value = "1.36"
$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0334"
$ws.Range("E50").Value = "  +7.14%  "

# Row 51 - OKB
$ws.Range("D51").Value = "43.93"
$ws.Range("E51").Value = "  -0.63%  "
